$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55; this shifts the existing rows 55-98
# down to 56-99 (and the sheet dimension grows to R99 automatically).
$ws.Rows("55").Insert()

# Populate the newly inserted row 55 with the new weekly record.
$ws.Range("A55").Value = 10
$ws.Range("B55").Value = "Vega Modelo de Temuco"
$ws.Range("C55").Value = "La Araucanía"
$ws.Range("D55").Value = 45040
$ws.Range("E55").Value = 9
$ws.Range("F55").Value = 300000001
$ws.Range("G55").Value = "Rabanito"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 40
$ws.Range("K55").Value = 8000
$ws.Range("L55").Value = 8000
$ws.Range("M55").Value = 8000
$ws.Range("N55").Value = "$/docena de paquetes"
$ws.Range("O55").Value = "Provincia de Cautín"
$ws.Range("P55").Value = 667
$ws.Range("Q55").Value = 12
$ws.Range("R55").Value = "Hortaliza"
